# "Avances de una semana"
#
# Original workbook has a single sheet ("Sheet1") containing a lookup table of
# electrode-pair labels, sorted alphabetically. This edit:
#   1. Renames Sheet1 -> "Alfabetico" (keeps its data/order/selection as-is,
#      just stops being the selected tab)
#   2. Adds three new sheets after it - "Stam", "Inter", "Redes" - which hold
#      the same label/filename pairs but re-ordered (grouped by channel
#      family instead of alphabetically) and carrying two extra numeric
#      columns ("x","y").
#   3. Makes "Redes" (the last sheet) the active/selected tab.

$wb = $excel.ActiveWorkbook
$alfa = $wb.Worksheets.Item(1)
$alfa.Name = "Alfabetico"

# ---------------------------------------------------------------------------
# New sheets, created right after "Alfabetico" in order so tab order comes
# out as Alfabetico, Stam, Inter, Redes.
# ---------------------------------------------------------------------------
$stam = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $alfa)
$stam.Name = "Stam"

$inter = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $stam)
$inter.Name = "Inter"

$redes = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $inter)
$redes.Name = "Redes"

# ---------------------------------------------------------------------------
# Helper data: label / filename pairs, in "natural" (non-alphabetical) order,
# as entered on "Stam"/"Inter". First 9 rows are single-channel pairs (which
# also carry x,y numbers on Stam/Inter); last 4 rows are triple/compound
# channel combos (which carry their own, different x,y numbers on Redes).
# ---------------------------------------------------------------------------
$singles = @(
    @("Fp1-Fp2", "FP1_FP2", 1, 5),
    @("F7-F8",   "F7_F8",   2, 4),
    @("F3-F4",   "F3_F4",   1, 4),
    @("T3-T4",   "T3_T4",   2, 3),
    @("C3-C4",   "C3_C4",   1, 3),
    @("T5-T6",   "T5_T6",   2, 2),
    @("P3-P4",   "P3_P4",   1, 2),
    @("O1-O2",   "O1_O2",   1, 1),
    @("LOG-ROG", "LOG_ROG", 1, 6)
)

$combos = @(
    @("Fp2-P4",    "FP2_P4",    2, 2),
    @("Fp1-P3",    "FP1_P3",    1, 2),
    @("O2-P4-T4",  "O2_P4_T4",  2, 1),
    @("O1-P3-T3",  "O1_P3_T3",  1, 1)
)

# ---------------------------------------------------------------------------
# Style helper: copy the number-format/fill/font of a reference cell on
# "Alfabetico" (styles 1/2/3 already exist in the workbook) onto a target
# range, using copy + paste-formats so we reuse the existing style indices
# instead of inventing new ones.
# ---------------------------------------------------------------------------
function Copy-Style($srcRange, $dstRange) {
    $srcRange.Copy()
    $dstRange.PasteSpecial(-4122)
}

$xlPasteFormats = -4122

# Reference style rows already present on Alfabetico:
#   row1 -> style 1 (bold header)
#   row6 (O1-O2) -> style 2
#   row3 (F3-F4) -> style 3
$styleHeaderRow = $alfa.Range("A1:B1")
$style2Row      = $alfa.Range("A6:B6")
$style3Row      = $alfa.Range("A3:B3")

# ===========================================================================
# Stam : header + 9 singles (with x,y) + 4 combos (no x,y)
# ===========================================================================
$stam.Range("A1").Value = "Etiqueta"
$stam.Range("B1").Value = "Nombre_archivo"
$stam.Range("C1").Value = "x"
$stam.Range("D1").Value = "y"
Copy-Style $styleHeaderRow $stam.Range("A1:B1")

$r = 2
foreach ($row in $singles) {
    $stam.Cells.Item($r, 1).Value = $row[0]
    $stam.Cells.Item($r, 2).Value = $row[1]
    $stam.Cells.Item($r, 3).Value = $row[2]
    $stam.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}
foreach ($row in $combos) {
    $stam.Cells.Item($r, 1).Value = $row[0]
    $stam.Cells.Item($r, 2).Value = $row[1]
    $r = $r + 1
}

# Row styles follow the same alternating 2/3/3/3 pattern (by groups of 3)
# that "Alfabetico" uses, in Stam's natural row order (rows 2..14):
$stamStyles = @(2, 2, 2, 3, 3, 3, 2, 2, 3, 2, 2, 3, 3)
for ($i = 0; $i -lt $stamStyles.Count; $i++) {
    $rowNum = $i + 2
    $target = $stam.Range("A" + $rowNum + ":B" + $rowNum)
    if ($stamStyles[$i] -eq 2) {
        Copy-Style $style2Row $target
    } else {
        Copy-Style $style3Row $target
    }
}

$stam.Columns.Item(2).AutoFit()
$stam.Range("A1:XFD1048576").Select()

# ===========================================================================
# Inter : header + the same 9 singles (with x,y) - no combos
# ===========================================================================
$inter.Range("A1").Value = "Etiqueta"
$inter.Range("B1").Value = "Nombre_archivo"
$inter.Range("C1").Value = "x"
$inter.Range("D1").Value = "y"
Copy-Style $styleHeaderRow $inter.Range("A1:B1")

$r = 2
foreach ($row in $singles) {
    $inter.Cells.Item($r, 1).Value = $row[0]
    $inter.Cells.Item($r, 2).Value = $row[1]
    $inter.Cells.Item($r, 3).Value = $row[2]
    $inter.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}

$interStyles = @(2, 2, 2, 3, 3, 3, 2, 2, 3)
for ($i = 0; $i -lt $interStyles.Count; $i++) {
    $rowNum = $i + 2
    $target = $inter.Range("A" + $rowNum + ":B" + $rowNum)
    if ($interStyles[$i] -eq 2) {
        Copy-Style $style2Row $target
    } else {
        Copy-Style $style3Row $target
    }
}

$inter.Columns.Item(2).AutoFit()
$inter.Range("F8").Select()

# ===========================================================================
# Redes : header + the 4 combos, each with its OWN x,y values
# ===========================================================================
$redes.Range("A1").Value = "Etiqueta"
$redes.Range("B1").Value = "Nombre_archivo"
$redes.Range("C1").Value = "x"
$redes.Range("D1").Value = "y"
Copy-Style $styleHeaderRow $redes.Range("A1:B1")

$r = 2
foreach ($row in $combos) {
    $redes.Cells.Item($r, 1).Value = $row[0]
    $redes.Cells.Item($r, 2).Value = $row[1]
    $redes.Cells.Item($r, 3).Value = $row[2]
    $redes.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}

$redesStyles = @(2, 2, 3, 3)
for ($i = 0; $i -lt $redesStyles.Count; $i++) {
    $rowNum = $i + 2
    $target = $redes.Range("A" + $rowNum + ":B" + $rowNum)
    if ($redesStyles[$i] -eq 2) {
        Copy-Style $style2Row $target
    } else {
        Copy-Style $style3Row $target
    }
}

$redes.Range("E5").Select()

# ---------------------------------------------------------------------------
# "Alfabetico" keeps its own prior selection; just make sure the tab
# selection marker moves to "Redes", the new active sheet.
# ---------------------------------------------------------------------------
$redes.Activate()
